$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Se agrega ambiente QA en creacion de Plan Autonomia (nueva fila 3)
$ws.Range("B3").Value = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/pc/PolicyCenter.do"
$ws.Range("A3").Value = "ssurgwsoadev4-oci.opc.oracleoutsourcing.com"
$ws.Range("C3").Value = "su"
$ws.Range("D3").Value = "gw"

# Hipervinculo de B3 apuntando a la nueva URL, con el mismo formato que B2
$ws.Hyperlinks.Add($ws.Range("B3"), "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/pc/PolicyCenter.do") | Out-Null
$ws.Range("B3").Style = $ws.Range("B2").Style

$ws.Range("D4").Select() | Out-Null
